# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps on the Overview, zh-cn, and de-de
# sheets to reflect the newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-25 05:03:14"

$wsZhCn.Range("H2").Value = "2016-08-25 05:03:08"
$wsZhCn.Range("K2").Value = "2016-08-25 05:03:25"

$wsDeDe.Range("H2").Value = "2016-08-25 05:03:14"
$wsDeDe.Range("K2").Value = "2016-08-25 05:03:33"
